# Actualizacion automatica del tracker
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Resolve the pending result for the existing row 84 ---
$ws.Range("G84").Value = "Fallo"
$ws.Range("H84").Value = -1

# --- Append newly tracked matches as rows 87-89 ---

# Row 87
$ws.Range("A87").Value = 14821167
$ws.Range("B87").NumberFormat = "@"
$ws.Range("B87").Value = "2025-10-04"
$ws.Range("B87").Style = "Normal"
$ws.Range("C87").Value = "Yafan Wang"
$ws.Range("D87").Value = "Zeynep Sonmez"
$ws.Range("E87").Value = "Gana Zeynep Sonmez"
$ws.Range("F87").Value = 1.83
$ws.Range("G87").Value = "'"
$ws.Range("G87").Style = "Normal"
$ws.Range("H87").Value = "'"
$ws.Range("H87").Style = "Normal"

# Row 88
$ws.Range("A88").Value = 14763098
$ws.Range("B88").NumberFormat = "@"
$ws.Range("B88").Value = "2025-10-04"
$ws.Range("B88").Style = "Normal"
$ws.Range("C88").Value = "Darja Vidmanova"
$ws.Range("D88").Value = "Nikola Bartunkova"
$ws.Range("E88").Value = "Gana Nikola Bartunkova"
$ws.Range("F88").Value = 1.83
$ws.Range("G88").Value = "'"
$ws.Range("G88").Style = "Normal"
$ws.Range("H88").Value = "'"
$ws.Range("H88").Style = "Normal"

# Row 89
$ws.Range("A89").Value = 14763655
$ws.Range("B89").NumberFormat = "@"
$ws.Range("B89").Value = "2025-10-04"
$ws.Range("B89").Style = "Normal"
$ws.Range("C89").Value = "Viktorija Golubic"
$ws.Range("D89").Value = "Tatjana Maria"
$ws.Range("E89").Value = "Gana Tatjana Maria"
$ws.Range("F89").Value = 2.25
$ws.Range("G89").Value = "'"
$ws.Range("G89").Style = "Normal"
$ws.Range("H89").Value = "'"
$ws.Range("H89").Style = "Normal"
